$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 89115
$ws.Range("E2").Value = 1618
$ws.Range("F2").Value = 1618
$ws.Range("G2").Value = 606
$ws.Range("H2").Value = 564
$ws.Range("I2").Value = 550
$ws.Range("J2").Value = 14
$ws.Range("K2").Value = 61519
$ws.Range("L2").Value = 51985
$ws.Range("M2").Value = 9533
$ws.Range("N2").Value = 9527
$ws.Range("O2").Value = 6
$ws.Range("P2").Value = 2000
$ws.Range("Q2").Value = 2777
$ws.Range("R2").Value = -2219
$ws.Range("S2").Value = 498
$ws.Range("T2").Value = 679
$ws.Range("U2").Value = 2098
$ws.Range("V2").Value = 15388
$ws.Range("W2").Value = 1.82
$ws.Range("X2").Value = 0.63
$ws.Range("Y2").Value = 5.92
$ws.Range("Z2").Value = 0.93
$ws.Range("AA2").Value = 545.29
$ws.Range("AB2").Value = 531.22
$ws.Range("AC2").Value = 857
$ws.Range("AD2").Value = 27.58
$ws.Range("AE2").Value = 15578
$ws.Range("AF2").Value = 1.52
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 64181675

# Row 3
$ws.Range("D3").Value = 64413
$ws.Range("E3").Value = -14543
$ws.Range("F3").Value = -14543
$ws.Range("G3").Value = -14501
$ws.Range("H3").Value = -13043
$ws.Range("I3").Value = -13053
$ws.Range("J3").Value = 9
$ws.Range("K3").Value = 56308
$ws.Range("L3").Value = 59437
$ws.Range("M3").Value = -3129
$ws.Range("N3").Value = -3145
$ws.Range("O3").Value = 16
$ws.Range("P3").Value = 2000
$ws.Range("Q3").Value = -8335
$ws.Range("R3").Value = 204
$ws.Range("S3").Value = 10724
$ws.Range("T3").Value = 202
$ws.Range("U3").Value = -8537
$ws.Range("V3").Value = 28221
$ws.Range("W3").Value = -22.58
$ws.Range("X3").Value = -20.25
$ws.Range("Y3").Value = -409.04
$ws.Range("Z3").Value = -22.14
$ws.Range("AA3").Value = -1899.34
$ws.Range("AB3").Value = -121.43
$ws.Range("AC3").Value = -20337
$ws.Range("AD3").Value = -0.72
$ws.Range("AE3").Value = -5142
$ws.Range("AF3").Value = -2.83
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 64181675

# Row 4
$ws.Range("D4").Value = 70094
$ws.Range("E4").Value = 701
$ws.Range("F4").Value = 701
$ws.Range("G4").Value = 564
$ws.Range("H4").Value = 94
$ws.Range("I4").Value = 258
$ws.Range("J4").Value = -163
$ws.Range("K4").Value = 55201
$ws.Range("L4").Value = 45239
$ws.Range("M4").Value = 9962
$ws.Range("N4").Value = 10109
$ws.Range("O4").Value = -147
$ws.Range("P4").Value = 9800
$ws.Range("Q4").Value = 1879
$ws.Range("R4").Value = 944
$ws.Range("S4").Value = -1019
$ws.Range("T4").Value = 241
$ws.Range("U4").Value = 1638
$ws.Range("V4").Value = 14516
$ws.Range("W4").Value = 1
$ws.Range("X4").Value = 0.13
$ws.Range("Y4").Value = 7.4
$ws.Range("Z4").Value = 0.17
$ws.Range("AA4").Value = 454.12
$ws.Range("AB4").Value = 26.19
$ws.Range("AC4").Value = 190
$ws.Range("AD4").Value = 54.29
$ws.Range("AE4").Value = 5158
$ws.Range("AF4").Value = 2
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 196000000
$ws.Range("AG4").ClearContents()
$ws.Range("AH4").ClearContents()

# Row 5
$ws.Range("D5").Value = 55362
$ws.Range("E5").Value = 469
$ws.Range("F5").Value = 469
$ws.Range("G5").Value = -296
$ws.Range("H5").Value = -521
$ws.Range("I5").Value = -454
$ws.Range("J5").Value = -68
$ws.Range("K5").Value = 50838
$ws.Range("L5").Value = 40792
$ws.Range("M5").Value = 10047
$ws.Range("N5").Value = 10262
$ws.Range("O5").Value = -215
$ws.Range("P5").Value = 9800
$ws.Range("Q5").Value = -1340
$ws.Range("R5").Value = -152
$ws.Range("S5").Value = 381
$ws.Range("T5").Value = 138
$ws.Range("U5").Value = -1478
$ws.Range("V5").Value = 14449
$ws.Range("W5").Value = 0.85
$ws.Range("X5").Value = -0.9399999999999999
$ws.Range("Y5").Value = -4.45
$ws.Range("Z5").Value = -0.98
$ws.Range("AA5").Value = 406.02
$ws.Range("AB5").Value = 15.33
$ws.Range("AC5").Value = -231
$ws.Range("AD5").Value = -53.59
$ws.Range("AE5").Value = 5236
$ws.Range("AF5").Value = 2.37
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 196000000
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()

# Row 6
$ws.Range("D6").Value = 54798
$ws.Range("E6").Value = 2061
$ws.Range("F6").Value = 2061
$ws.Range("G6").Value = 1702
$ws.Range("H6").Value = 702
$ws.Range("I6").Value = 687
$ws.Range("K6").Value = 46285
$ws.Range("L6").Value = 35947
$ws.Range("M6").Value = 10338
$ws.Range("N6").Value = 10550
$ws.Range("P6").Value = 9800
$ws.Range("Q6").Value = 4462
$ws.Range("R6").Value = 207
$ws.Range("S6").Value = -8492
$ws.Range("T6").Value = 146
$ws.Range("U6").Value = 4316
$ws.Range("V6").Value = 5575
$ws.Range("W6").Value = 3.76
$ws.Range("X6").Value = 1.28
$ws.Range("Y6").Value = 6.6
$ws.Range("Z6").Value = 1.45
$ws.Range("AA6").Value = 347.73
$ws.Range("AB6").Value = 21.87
$ws.Range("AC6").Value = 350
$ws.Range("AD6").Value = 50.22
$ws.Range("AE6").Value = 5383
$ws.Range("AF6").Value = 3.27
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 196000000
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Row 7
$ws.Range("D7").Value = 62711
$ws.Range("E7").Value = 4079
$ws.Range("G7").Value = 4263
$ws.Range("H7").Value = 3202
$ws.Range("I7").Value = 3190
$ws.Range("K7").Value = 47433
$ws.Range("L7").Value = 33887
$ws.Range("M7").Value = 13547
$ws.Range("N7").Value = 13769
$ws.Range("P7").Value = 9800
$ws.Range("Q7").Value = 4079
$ws.Range("R7").Value = -75
$ws.Range("S7").Value = -2370
$ws.Range("T7").Value = 126
$ws.Range("U7").Value = 4051
$ws.Range("W7").Value = 6.5
$ws.Range("X7").Value = 5.11
$ws.Range("Y7").Value = 26.23
$ws.Range("Z7").Value = 6.83
$ws.Range("AA7").Value = 250.15
$ws.Range("AC7").Value = 1627
$ws.Range("AD7").Value = 11.24
$ws.Range("AE7").Value = 7025
$ws.Range("AF7").Value = 2.6
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 66530
$ws.Range("E8").Value = 4010
$ws.Range("G8").Value = 4483
$ws.Range("H8").Value = 3310
$ws.Range("I8").Value = 3272
$ws.Range("K8").Value = 51205
$ws.Range("L8").Value = 34367
$ws.Range("M8").Value = 16837
$ws.Range("N8").Value = 17061
$ws.Range("P8").Value = 9800
$ws.Range("Q8").Value = 2775
$ws.Range("R8").Value = -170
$ws.Range("S8").Value = -542
$ws.Range("T8").Value = 171
$ws.Range("U8").Value = 2576
$ws.Range("W8").Value = 6.03
$ws.Range("X8").Value = 4.98
$ws.Range("Y8").Value = 21.22
$ws.Range("Z8").Value = 6.72
$ws.Range("AA8").Value = 204.12
$ws.Range("AC8").Value = 1669
$ws.Range("AD8").Value = 10.33
$ws.Range("AE8").Value = 8705
$ws.Range("AF8").Value = 1.98
$ws.Range("AG8").Value = 26
$ws.Range("AH8").Value = 0.15
$ws.Range("AI8").Value = 1.59

# Row 9
$ws.Range("D9").Value = 72142
$ws.Range("E9").Value = 4588
$ws.Range("G9").Value = 5053
$ws.Range("H9").Value = 3735
$ws.Range("I9").Value = 3757
$ws.Range("K9").Value = 55215
$ws.Range("L9").Value = 34745
$ws.Range("M9").Value = 20470
$ws.Range("N9").Value = 20740
$ws.Range("P9").Value = 9800
$ws.Range("Q9").Value = 3527
$ws.Range("R9").Value = -171
$ws.Range("S9").Value = -455
$ws.Range("T9").Value = 144
$ws.Range("U9").Value = 3324
$ws.Range("W9").Value = 6.36
$ws.Range("X9").Value = 5.18
$ws.Range("Y9").Value = 19.88
$ws.Range("Z9").Value = 7.02
$ws.Range("AA9").Value = 169.74
$ws.Range("AC9").Value = 1917
$ws.Range("AD9").Value = 9
$ws.Range("AE9").Value = 10582
$ws.Range("AF9").Value = 1.63
$ws.Range("AG9").Value = 56
$ws.Range("AH9").Value = 0.33
$ws.Range("AI9").Value = 2.94
